$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (the "Förändrad" date column) for rows 2 through 233 all
# held the serial date value 45204 (2023-10-05) and were updated to
# 45205 (2023-10-06).
$ws.Range("C2:C233").Value = 45205
